$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: clear the now-unused "14"/"15" index markers (O1:P1). The rest of row 1
# (A1:N1 and Q1:Z1) is untouched.
$ws.Range("O1:P1").ClearContents()

# Row 2 (sensor/label row): one duplicate "Index + Middle" gesture column was
# removed, so C2:N2 now hold what used to live in D2:O2, and the trailing
# columns (O2:P2, including the removed "Thumb Down" label) are cleared.
$row2 = @("Open Hand ","Closed Hand","Index","Index + Middle","Index + Middle + Ring","Index + Middle + Ring + Pinky","Open Hand","Closed Hand","Index","Pinky","Index + Pinky","Index + Middle + Ring + Pinky","Thumb Left","Thumbs Up")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}
$ws.Range("O2:P2").ClearContents()

# Row 3 (gesture command row): same leftward shift starting at column C.
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "Left 90 degrees"
$ws.Range("J3").Value = "Right 90 degrees"
$ws.Range("K3").Value = "Forward"
$ws.Range("L3").Value = "Backwards"
$ws.Range("M3").Value = "Rotate 180 degrees"
$ws.Range("N3").Value = "Get Gestures"
$ws.Range("O3").ClearContents()

# Row 4 (hand row): same leftward shift starting at column C.
$row4 = @("Right","Right","Left","Left","Left","Left","Left","Left","Right","Right","Right","Right","Right","Right")
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}
$ws.Range("O4:P4").ClearContents()

# Restore the selection to match the author's final cursor position.
$ws.Range("J1").Select()
